# Updated BOM to reflect layout changes.
# Remove the "JP2Q / JP1" jumper-header line item (row 2) from the BOM table.
# Deleting the entire row shifts every row below it up by one, which is
# exactly what the target diff shows (old row N+1 == new row N for the rest
# of the table), and shrinks the Table1 / autofilter range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Select row 2 first (mirrors the manual "click row header, Delete Row" flow)
# and reposition the view like the saved file shows.
$ws.Activate()
$ws.Range("A2:XFD2").Select()

# Delete the entire row 2 (JP2Q / JP1 jumper header), shifting rows 3-30 up.
$ws.Rows("2:2").Delete()

# The 150uF / CAP_POLE (C26) line's quantity was updated from 1 to 32 as part
# of this layout revision (now row 12 after the shift).
$ws.Range("A12").Value = 32

# Restore the view/selection to match the post-edit state.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C1").Select()
$ws.Range("C2").Select()
